# This script edits the "Test Report" bug table:
#  - Row for bug 32 (Word table row 33): the "Scenario" cell text is re-typed,
#    which causes Word's grammar checker to wrap the word "wrong" with
#    gramStart/gramEnd proofing marks and merges the trailing ". " into the
#    sentence (instead of a separate bookmark + ". " run).
#  - Six new bug rows (33-37) are filled in with their scenario / solution /
#    related-bug text (these rows already existed as blank rows at the end
#    of the table, ready to be completed).
#
# All inserted paragraphs use InsertXML so that the run / proofErr structure
# matches exactly, and keep the existing "sz=18 / szCs=18" (9pt) formatting
# used throughout the table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Make-Para([string]$innerXml) {
    return '<w:p ' + $wNs + '><w:pPr><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr>' + $innerXml + '</w:p>'
}

function Make-Run([string]$text, [bool]$preserve) {
    if ($preserve) {
        $t1 = '<w:t xml:space="preserve">' + $text + '</w:t>'
    } else {
        $t1 = '<w:t>' + $text + '</w:t>'
    }
    return '<w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>' + $t1 + '</w:r>'
}

# ---------------------------------------------------------------------
# Row for bug 32 (Word row 33), column 2 (Scenario)
# NOTE: InsertXML always rewrites the whole paragraph it targets, so the
# complete paragraph (including the leading runs that are not touched by
# the diff) is reconstructed here; only the trailing sentence changes.
# ---------------------------------------------------------------------
$inner = ''
$inner += Make-Run 'As a customer, I ' $true
$inner += '<w:proofErr w:type="gramStart"/>'
$inner += Make-Run 'am able ' $true
$inner += Make-Run 'to' $false
$inner += '<w:proofErr w:type="gramEnd"/>'
$inner += Make-Run ' checkout with nothing in my cart. This doesn’t do anything ' $true
$inner += '<w:proofErr w:type="gramStart"/>'
$inner += Make-Run 'wrong' $false
$inner += '<w:proofErr w:type="gramEnd"/>'
$inner += Make-Run ' but this is not sensible. ' $true
$t.Cell(33, 2).Range.InsertXML((Make-Para $inner))

# ---------------------------------------------------------------------
# New row: bug 33 (Word row 34)
# ---------------------------------------------------------------------
$t.Cell(34, 1).Range.InsertXML((Make-Para (Make-Run '33' $false)))

$inner = ''
$inner += Make-Run 'When an admin ' $true
$inner += '<w:proofErr w:type="gramStart"/>'
$inner += Make-Run 'create' $false
$inner += '<w:proofErr w:type="gramEnd"/>'
$inner += Make-Run ' another admin through the add/edit user page, when the admin created, you receive a message saying that a customer is created (and not admin).' $true
$t.Cell(34, 2).Range.InsertXML((Make-Para $inner))

$t.Cell(34, 3).Range.InsertXML((Make-Para (Make-Run 'Change the message so that it says admin, when an admin is created. Check bug 34.' $false)))

$t.Cell(34, 4).Range.InsertXML((Make-Para (Make-Run '34' $false)))

# ---------------------------------------------------------------------
# New row: bug 34 (Word row 35)
# ---------------------------------------------------------------------
$t.Cell(35, 1).Range.InsertXML((Make-Para (Make-Run '34' $false)))

$inner = ''
$inner += Make-Run 'The roles feature does not work on the ' $true
$inner += '<w:proofErr w:type="gramStart"/>'
$inner += Make-Run 'admin’s' $false
$inner += '<w:proofErr w:type="gramEnd"/>'
$inner += Make-Run ' create user page. It will always create a customer' $true
$t.Cell(35, 2).Range.InsertXML((Make-Para $inner))

$t.Cell(35, 3).Range.InsertXML((Make-Para (Make-Run 'Make it so that an admin is created when the role is set to admin.' $false)))

# ---------------------------------------------------------------------
# New row: bug 35 (Word row 36)
# ---------------------------------------------------------------------
$t.Cell(36, 1).Range.InsertXML((Make-Para (Make-Run '35' $false)))

$t.Cell(36, 2).Range.InsertXML((Make-Para (Make-Run 'After I add an account to any Id, it takes me to a page where the buttons do not work. ' $true)))

$t.Cell(36, 3).Range.InsertXML((Make-Para (Make-Run 'Make it so that the after I create the account, it takes me to a working page with buttons.' $false)))

# ---------------------------------------------------------------------
# New row: bug 36 (Word row 37)
# ---------------------------------------------------------------------
$t.Cell(37, 1).Range.InsertXML((Make-Para (Make-Run '36' $false)))

$t.Cell(37, 2).Range.InsertXML((Make-Para (Make-Run 'There should be a notification saying that adding the account to the customer was successful (or unsuccessful if I try to add it to an admin).' $false)))

$t.Cell(37, 3).Range.InsertXML((Make-Para (Make-Run 'Make a message appear to show that it works' $false)))

# ---------------------------------------------------------------------
# New row: bug 37 (Word row 38)
# ---------------------------------------------------------------------
$t.Cell(38, 1).Range.InsertXML((Make-Para (Make-Run '37' $false)))

$t.Cell(38, 2).Range.InsertXML((Make-Para (Make-Run 'It is currently possible to add accounts to admin Ids.' $false)))

$t.Cell(38, 3).Range.InsertXML((Make-Para (Make-Run 'Make sure that the app prevents admins from adding accounts to their ids.' $false)))

$inner = Make-Run '36' $false
$inner += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$t.Cell(38, 4).Range.InsertXML((Make-Para $inner))

Write-Host "Edit complete"
